$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 999.4706
$ws.Range("I98").Value = 1027.1428
$ws.Range("K98").Value = 1027.1428
$ws.Range("M98").Value = 470.8571999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 999.4706
$ws.Range("I122").Value = 1027.1428
$ws.Range("K122").Value = 3081.4284
$ws.Range("M122").Value = -631.4284000000002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 3176
$ws.Range("I125").Value = 2669.3333
$ws.Range("J125").Value = 4392
$ws.Range("K125").Value = 24023.9997
$ws.Range("L125").Value = 39528
$ws.Range("M125").Value = -21563.9997
$ws.Range("N125").Value = -44448

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 9949.5
$ws.Range("J9").Value = 9949.5
$ws.Range("L9").Value = 9949.5
$ws.Range("N9").Value = -10289.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H20").Value = 9949.5
$ws.Range("J20").Value = 9949.5
$ws.Range("L20").Value = 9949.5
$ws.Range("N20").Value = -10489.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2308.3333
$ws.Range("I63").Value = 1645.4546
$ws.Range("J63").Value = 3350
$ws.Range("K63").Value = 1645.4546
$ws.Range("L63").Value = 3350
$ws.Range("M63").Value = -959.4546
$ws.Range("N63").Value = -4722

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2308.3333
$ws.Range("I66").Value = 1645.4546
$ws.Range("J66").Value = 3350
$ws.Range("K66").Value = 8227.273000000001
$ws.Range("L66").Value = 16750
$ws.Range("M66").Value = -4795.273000000001
$ws.Range("N66").Value = -23614

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 7862.0713
$ws.Range("J51").Value = 7844.5386
$ws.Range("L51").Value = 7844.5386
$ws.Range("N51").Value = -9316.5386

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 7862.0713
$ws.Range("J61").Value = 7844.5386
$ws.Range("L61").Value = 7844.5386
$ws.Range("N61").Value = -8540.5386

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 6946805.5
$ws.Range("I62").Value = 55555556
$ws.Range("J62").Value = 2698.2856
$ws.Range("K62").Value = 55555556
$ws.Range("L62").Value = 2698.2856
$ws.Range("M62").Value = -55554932
$ws.Range("N62").Value = -3946.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 6946805.5
$ws.Range("I65").Value = 55555556
$ws.Range("J65").Value = 2698.2856
$ws.Range("K65").Value = 277777780
$ws.Range("L65").Value = 13491.428
$ws.Range("M65").Value = -277774660
$ws.Range("N65").Value = -19731.428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 13442.9
$ws.Range("I99").Value = 4233.5
$ws.Range("J99").Value = 27257
$ws.Range("K99").Value = 4233.5
$ws.Range("L99").Value = 27257
$ws.Range("M99").Value = -2735.5
$ws.Range("N99").Value = -30253

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 941
$ws.Range("I122").Value = 866.6667
$ws.Range("J122").Value = 978.1667
$ws.Range("K122").Value = 2600.0001
$ws.Range("L122").Value = 2934.5001
$ws.Range("M122").Value = -150.0001000000002
$ws.Range("N122").Value = -7834.5001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 13442.9
$ws.Range("I126").Value = 4233.5
$ws.Range("J126").Value = 27257
$ws.Range("K126").Value = 12700.5
$ws.Range("L126").Value = 81771
$ws.Range("M126").Value = -10230.5
$ws.Range("N126").Value = -86711

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H129").Value = 49332.668
$ws.Range("J129").Value = 49332.668
$ws.Range("L129").Value = 49332.668
$ws.Range("N129").Value = -59332.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 34892.75
$ws.Range("J37").Value = 34892.75
$ws.Range("L37").Value = 104678.25
$ws.Range("N37").Value = -104902.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 356899
$ws.Range("I107").Value = 517.0833
$ws.Range("J107").Value = 515290.97
$ws.Range("K107").Value = 1551.2499
$ws.Range("L107").Value = 1545872.91
$ws.Range("M107").Value = 368.7501
$ws.Range("N107").Value = -1549712.91

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 849.23
$ws.Range("I131").Value = 630
$ws.Range("J131").Value = 860.76843
$ws.Range("K131").Value = 1890
$ws.Range("L131").Value = 2582.30529
$ws.Range("M131").Value = 3150
$ws.Range("N131").Value = -12662.30529

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3364.3333
$ws.Range("I102").Value = 3299.111
$ws.Range("J102").Value = 3560
$ws.Range("K102").Value = 3299.111
$ws.Range("L102").Value = 3560
$ws.Range("M102").Value = -1677.111
$ws.Range("N102").Value = -6804

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 934.6429000000001
$ws.Range("I122").Value = 907.9091
$ws.Range("J122").Value = 1032.6666
$ws.Range("K122").Value = 2723.7273
$ws.Range("L122").Value = 3097.9998
$ws.Range("M122").Value = -273.7273
$ws.Range("N122").Value = -7997.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3287.8635
$ws.Range("I126").Value = 2878.4614
$ws.Range("J126").Value = 3879.2222
$ws.Range("K126").Value = 8635.3842
$ws.Range("L126").Value = 11637.6666
$ws.Range("M126").Value = -6165.3842
$ws.Range("N126").Value = -16577.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3233.6
$ws.Range("I7").Value = 1937.375
$ws.Range("J7").Value = 4715
$ws.Range("K7").Value = 1937.375
$ws.Range("L7").Value = 4715
$ws.Range("M7").Value = -1825.375
$ws.Range("N7").Value = -4939

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 37450
$ws.Range("J36").Value = 37450
$ws.Range("L36").Value = 37450
$ws.Range("N36").Value = -38574

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2079.8
$ws.Range("I40").Value = 1324.75
$ws.Range("J40").Value = 2354.3635
$ws.Range("K40").Value = 1324.75
$ws.Range("L40").Value = 2354.3635
$ws.Range("M40").Value = -1188.75
$ws.Range("N40").Value = -2626.3635

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4633.3335
$ws.Range("I122").Value = 7333.3335
$ws.Range("J122").Value = 1933.3334
$ws.Range("K122").Value = 22000.0005
$ws.Range("L122").Value = 5800.0002
$ws.Range("M122").Value = -19550.0005
$ws.Range("N122").Value = -10700.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3233.6
$ws.Range("I126").Value = 1937.375
$ws.Range("J126").Value = 4715
$ws.Range("K126").Value = 5812.125
$ws.Range("L126").Value = 14145
$ws.Range("M126").Value = -3342.125
$ws.Range("N126").Value = -19085

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1800.0667
$ws.Range("I122").Value = 999.8570999999999
$ws.Range("K122").Value = 2999.5713
$ws.Range("M122").Value = -549.5712999999996

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1162.3928
$ws.Range("I126").Value = 1274.619
$ws.Range("J126").Value = 825.7143
$ws.Range("K126").Value = 3823.857
$ws.Range("L126").Value = 2477.1429
$ws.Range("M126").Value = -1353.857
$ws.Range("N126").Value = -7417.1429
